$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo/placeholder task description in row 7 (column A)
$ws.Range("A7").Value = "If we got our desired output, we start with implementing the view with the graphical output"

# Fill in progress ("Remain" worked so far) values in column D
$ws.Range("D3").Value = 10
$ws.Range("D6").Value = 18
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 0

# Row 7's "Curr. Est." (C7) becomes a fixed re-estimate instead of the shared formula
$ws.Range("C7").Value = 10

# Update the active selection to A17
$ws.Range("A17").Select()
